# Add the new "2022-Q1" quarterly sheet, inserted right before the "总计" summary sheet.
$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$ws = $wb.Worksheets.Add($totalSheet)
$ws.Name = "2022-Q1"

# Mirror the header/index-column formatting used by the other quarterly sheets.
$fmtSrc = $wb.Worksheets.Item("2021-Q4")
$fmtSrc.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$fmtSrc.Range("A2").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "'000727"
$ws.Range("B2").ClearFormats()
$ws.Range("C2").Value = "'融通健康产业灵活配置混合A"
$ws.Range("C2").ClearFormats()
$ws.Range("D2").Value = "'15.30"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'94.68"
$ws.Range("E2").ClearFormats()
$ws.Range("F2").Value = "'9.10"
$ws.Range("F2").ClearFormats()
$ws.Range("G2").Value = "'1.3923"
$ws.Range("G2").ClearFormats()
$ws.Range("H2").Value = 3

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "'006218"
$ws.Range("B3").ClearFormats()
$ws.Range("C3").Value = "'富国生物医药科技混合A"
$ws.Range("C3").ClearFormats()
$ws.Range("D3").Value = "'9.55"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'87.25"
$ws.Range("E3").ClearFormats()
$ws.Range("F3").Value = "'7.77"
$ws.Range("F3").ClearFormats()
$ws.Range("G3").Value = "'0.7420"
$ws.Range("G3").ClearFormats()
$ws.Range("H3").Value = 2

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "'100016"
$ws.Range("B4").ClearFormats()
$ws.Range("C4").Value = "'富国天源沪港深平衡混合"
$ws.Range("C4").ClearFormats()
$ws.Range("D4").Value = "'6.23"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'72.29"
$ws.Range("E4").ClearFormats()
$ws.Range("F4").Value = "'7.01"
$ws.Range("F4").ClearFormats()
$ws.Range("G4").Value = "'0.4367"
$ws.Range("G4").ClearFormats()
$ws.Range("H4").Value = 1

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "'009274"
$ws.Range("B5").ClearFormats()
$ws.Range("C5").Value = "'融通健康产业灵活配置混合C"
$ws.Range("C5").ClearFormats()
$ws.Range("D5").Value = "'3.16"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'94.68"
$ws.Range("E5").ClearFormats()
$ws.Range("F5").Value = "'9.10"
$ws.Range("F5").ClearFormats()
$ws.Range("G5").Value = "'0.2876"
$ws.Range("G5").ClearFormats()
$ws.Range("H5").Value = 3

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "'001048"
$ws.Range("B6").ClearFormats()
$ws.Range("C6").Value = "'富国新兴产业股票"
$ws.Range("C6").ClearFormats()
$ws.Range("D6").Value = "'7.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'92.70"
$ws.Range("E6").ClearFormats()
$ws.Range("F6").Value = "'2.78"
$ws.Range("F6").ClearFormats()
$ws.Range("G6").Value = "'0.2221"
$ws.Range("G6").ClearFormats()
$ws.Range("H6").Value = 9

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "'011308"
$ws.Range("B7").ClearFormats()
$ws.Range("C7").Value = "'富国生物医药科技混合C"
$ws.Range("C7").ClearFormats()
$ws.Range("D7").Value = "'1.26"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'87.25"
$ws.Range("E7").ClearFormats()
$ws.Range("F7").Value = "'7.77"
$ws.Range("F7").ClearFormats()
$ws.Range("G7").Value = "'0.0979"
$ws.Range("G7").ClearFormats()
$ws.Range("H7").Value = 2

# Update the "总计" summary sheet with the new 2022-Q1 row, shifting existing rows down
# and renumbering the index column.
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 6
$total.Range("D2").Value = 3.18

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
